# Auto-generated Excel COM-interop script to apply Mandragora_Profits.xlsx edits
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 69547.39
$ws.Range("I86").Value = 95065.62
$ws.Range("K86").Value = 95065.62
$ws.Range("M86").Value = -93942.62
$ws.Range("H87").Value = 39924.5
$ws.Range("J87").Value = 39924.5
$ws.Range("L87").Value = 39924.5
$ws.Range("N87").Value = -42420.5
$ws.Range("H89").Value = 69547.39
$ws.Range("I89").Value = 95065.62
$ws.Range("K89").Value = 475328.1
$ws.Range("M89").Value = -469712.1
$ws.Range("H90").Value = 39924.5
$ws.Range("J90").Value = 39924.5
$ws.Range("L90").Value = 119773.5
$ws.Range("N90").Value = -132253.5
$ws.Range("H116").Value = 1162129.5
$ws.Range("I116").Value = 3061.4443
$ws.Range("J116").Value = 2652359.8
$ws.Range("K116").Value = 3061.4443
$ws.Range("L116").Value = 2652359.8
$ws.Range("M116").Value = 380.5556999999999
$ws.Range("N116").Value = -2659243.8
$ws.Range("H132").Value = 4027.125
$ws.Range("I132").Value = 3266.3403
$ws.Range("J132").Value = 5457.4
$ws.Range("K132").Value = 9799.0209
$ws.Range("L132").Value = 16372.2
$ws.Range("M132").Value = -7269.0209
$ws.Range("N132").Value = -21432.2
$ws.Range("H137").Value = 1854.619
$ws.Range("I137").Value = 1804.5927
$ws.Range("J137").Value = 1944.6666
$ws.Range("K137").Value = 5413.7781
$ws.Range("L137").Value = 5833.9998
$ws.Range("M137").Value = -2863.7781
$ws.Range("N137").Value = -10933.9998

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1558.381
$ws.Range("I74").Value = 1265.125
$ws.Range("J74").Value = 1949.3889
$ws.Range("K74").Value = 1265.125
$ws.Range("L74").Value = 1949.3889
$ws.Range("M74").Value = -391.125
$ws.Range("N74").Value = -3697.3889
$ws.Range("H77").Value = 1558.381
$ws.Range("I77").Value = 1265.125
$ws.Range("J77").Value = 1949.3889
$ws.Range("K77").Value = 6325.625
$ws.Range("L77").Value = 9746.9445
$ws.Range("M77").Value = -1957.625
$ws.Range("N77").Value = -18482.9445
$ws.Range("H110").Value = 2473.56
$ws.Range("I110").Value = 2157.6155
$ws.Range("J110").Value = 2815.8333
$ws.Range("K110").Value = 2157.6155
$ws.Range("L110").Value = 2815.8333
$ws.Range("M110").Value = -112.6154999999999
$ws.Range("N110").Value = -6905.8333
$ws.Range("H132").Value = 4249.857
$ws.Range("I132").Value = 1976.0667
$ws.Range("J132").Value = 7840.0527
$ws.Range("K132").Value = 5928.2001
$ws.Range("L132").Value = 23520.1581
$ws.Range("M132").Value = -3398.2001
$ws.Range("N132").Value = -28580.1581

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2000
$ws.Range("I105").Value = 1800
$ws.Range("J105").Value = 2600
$ws.Range("K105").Value = 1800
$ws.Range("L105").Value = 2600
$ws.Range("M105").Value = -53
$ws.Range("N105").Value = -6094
$ws.Range("H107").Value = 1381.3529
$ws.Range("I107").Value = 1463.0714
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 1463.0714
$ws.Range("L107").Value = 1000
$ws.Range("M107").Value = 456.9286
$ws.Range("N107").Value = -4840

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3380.0908
$ws.Range("I16").Value = 3523.875
$ws.Range("J16").Value = 2996.6667
$ws.Range("K16").Value = 3523.875
$ws.Range("L16").Value = 2996.6667
$ws.Range("M16").Value = -3236.875
$ws.Range("N16").Value = -3570.6667
$ws.Range("H31").Value = 7938810
$ws.Range("I31").Value = 1739.6842
$ws.Range("J31").Value = 14495520
$ws.Range("K31").Value = 1739.6842
$ws.Range("L31").Value = 14495520
$ws.Range("M31").Value = -1444.6842
$ws.Range("N31").Value = -14496110
$ws.Range("H34").Value = 7938810
$ws.Range("I34").Value = 1739.6842
$ws.Range("J34").Value = 14495520
$ws.Range("K34").Value = 1739.6842
$ws.Range("L34").Value = 14495520
$ws.Range("M34").Value = -1537.6842
$ws.Range("N34").Value = -14495924
$ws.Range("H41").Value = 17500
$ws.Range("I41").Value = 10000
$ws.Range("K41").Value = 10000
$ws.Range("M41").Value = -9572
$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").ClearContents()
$ws.Range("H51").Value = 3000
$ws.Range("I51").Value = 3000
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 3000
$ws.Range("L51").Value = 0
$ws.Range("M51").Value = -2264
$ws.Range("N51").ClearContents()
$ws.Range("H60").Value = 10218.6
$ws.Range("I60").Value = 6546.5
$ws.Range("J60").Value = 12666.667
$ws.Range("K60").Value = 6546.5
$ws.Range("L60").Value = 12666.667
$ws.Range("M60").Value = -6035.5
$ws.Range("N60").Value = -13688.667
$ws.Range("H61").Value = 3000
$ws.Range("I61").Value = 3000
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 3000
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -2652
$ws.Range("N61").ClearContents()
$ws.Range("H74").Value = 33552.75
$ws.Range("J74").Value = 33552.75
$ws.Range("L74").Value = 33552.75
$ws.Range("N74").Value = -35300.75
$ws.Range("H77").Value = 33552.75
$ws.Range("J77").Value = 33552.75
$ws.Range("L77").Value = 100658.25
$ws.Range("N77").Value = -109394.25
$ws.Range("H113").Value = 3380.0908
$ws.Range("I113").Value = 3523.875
$ws.Range("J113").Value = 2996.6667
$ws.Range("K113").Value = 3523.875
$ws.Range("L113").Value = 2996.6667
$ws.Range("M113").Value = -1353.875
$ws.Range("N113").Value = -7336.6667
$ws.Range("H132").Value = 2950.853
$ws.Range("I132").Value = 2153.2354
$ws.Range("J132").Value = 3748.4707
$ws.Range("K132").Value = 6459.706200000001
$ws.Range("L132").Value = 11245.4121
$ws.Range("M132").Value = -3929.706200000001
$ws.Range("N132").Value = -16305.4121

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 982.63464
$ws.Range("I68").Value = 904.7619
$ws.Range("J68").Value = 1035.3871
$ws.Range("K68").Value = 2714.2857
$ws.Range("L68").Value = 3106.1613
$ws.Range("M68").Value = -1903.2857
$ws.Range("N68").Value = -4728.1613
$ws.Range("H71").Value = 982.63464
$ws.Range("I71").Value = 904.7619
$ws.Range("J71").Value = 1035.3871
$ws.Range("K71").Value = 8142.857099999999
$ws.Range("L71").Value = 9318.483899999999
$ws.Range("M71").Value = -4086.857099999999
$ws.Range("N71").Value = -17430.4839

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3220.7896
$ws.Range("I80").Value = 2658.182
$ws.Range("J80").Value = 3994.375
$ws.Range("K80").Value = 2658.182
$ws.Range("L80").Value = 3994.375
$ws.Range("M80").Value = -1660.182
$ws.Range("N80").Value = -5990.375
$ws.Range("H83").Value = 3220.7896
$ws.Range("I83").Value = 2658.182
$ws.Range("J83").Value = 3994.375
$ws.Range("K83").Value = 13290.91
$ws.Range("L83").Value = 19971.875
$ws.Range("M83").Value = -8298.91
$ws.Range("N83").Value = -29955.875
$ws.Range("H132").Value = 3209106.2
$ws.Range("I132").Value = 5210936
$ws.Range("J132").Value = 6178.4
$ws.Range("K132").Value = 15632808
$ws.Range("L132").Value = 18535.2
$ws.Range("M132").Value = -15630278
$ws.Range("N132").Value = -23595.2
$ws.Range("H141").Value = 29000
$ws.Range("J141").Value = 29000
$ws.Range("L141").Value = 29000
$ws.Range("N141").Value = -39360

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1949
$ws.Range("I7").Value = 1520.909
$ws.Range("K7").Value = 1520.909
$ws.Range("M7").Value = -1408.909
$ws.Range("H122").Value = 4167.1113
$ws.Range("I122").Value = 3438
$ws.Range("J122").Value = 10000
$ws.Range("K122").Value = 10314
$ws.Range("L122").Value = 30000
$ws.Range("M122").Value = -7864
$ws.Range("N122").Value = -34900
$ws.Range("H126").Value = 1949
$ws.Range("I126").Value = 1520.909
$ws.Range("K126").Value = 4562.727000000001
$ws.Range("M126").Value = -2092.727000000001

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 4647.057
$ws.Range("I122").Value = 3659.6365
$ws.Range("J122").Value = 6318.077
$ws.Range("K122").Value = 10978.9095
$ws.Range("L122").Value = 18954.231
$ws.Range("M122").Value = -8528.9095
$ws.Range("N122").Value = -23854.231
$ws.Range("H132").Value = 2081.653
$ws.Range("I132").Value = 1685.1
$ws.Range("J132").Value = 3844.111
$ws.Range("K132").Value = 5055.299999999999
$ws.Range("L132").Value = 11532.333
$ws.Range("M132").Value = -2525.299999999999
$ws.Range("N132").Value = -16592.333
$ws.Range("H140").Value = 38885.8
$ws.Range("J140").Value = 38885.8
$ws.Range("L140").Value = 38885.8
$ws.Range("N140").Value = -49245.8
